# Automatische test-sync: 2025-06-20 08:30:50
# Adds a new incoming mail-log row (row 3) to the "Logs" sheet, updates the
# matching "Dashboard" pivot-style summary row, and extends the chart series
# / conditional-formatting ranges so they cover the new row as well.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append the new mail entry on row 3
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Is product X op voorraad?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Range("D3").Value = "Productinformatie"
$logs.Range("F3").Value = "2025-06-20 08:30:14"
$logs.Range("G3").Value = "Nee"

# Extend the existing conditional formatting (Categorie / Beantwoord columns)
# so it also covers the freshly added row.
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))

# ---------------------------------------------------------------------------
# 2. Dashboard sheet: add the matching category/count row
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 1

# ---------------------------------------------------------------------------
# 3. Chart: extend the category/value series so they include the new row
# ---------------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
